# The visible text "<id>p123r_1</id>" is split across three runs in the
# source paragraph: one run for the "<id>" tag (Courier New / 7f6000 /
# 9pt), a plain run for "p123r_1", and another Courier New / 7f6000 / 9pt
# run for "</id>". Re-run a Find & Replace over that exact span with the
# identical text so Word collapses the three runs into a single run,
# taking on the formatting of the first run in the found range.
$d = $word.ActiveDocument
$d.Content.Find.Execute("<id>p123r_1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p123r_1</id>", 2)
